$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("B4").Value = 597
$ws.Range("C4").Value = 572
$ws.Range("D4").Value = 567
$ws.Range("E4").Value = 587
$ws.Range("F4").Value = 576
$ws.Range("G4").Value = 563
$ws.Range("H4").Value = 562
$ws.Range("I4").Value = 571

$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("B5").Value = 671
$ws.Range("C5").Value = 649
$ws.Range("D5").Value = 639
$ws.Range("E5").Value = 655
$ws.Range("F5").Value = 639
$ws.Range("G5").Value = 622
$ws.Range("H5").Value = 617
$ws.Range("I5").Value = 632

$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Tsalenjikha Municipality"
Write-Output "ok"
